$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells receiving numeric-looking text must be forced to Text format
# so Excel does not silently convert them to numbers and lose exact formatting.
$textRefs = @("D5","D6","D8","D9","D10","D11","D14","D15","D16","D17","D20","D22","D23","D26","D27","D29","D30","D31","D32","D33","D34","D36","D38","D39","D40","D41","D44","D45","D46","D47","D48","D49","D51")
foreach ($r in $textRefs) {
    $ws.Range($r).NumberFormat = "@"
}

$ws.Range("D2").Value = "36.601.97"
$ws.Range("E2").Value = "  -1.44%  "
$ws.Range("D3").Value = "2.031.03"
$ws.Range("E3").Value = "  +0.94%  "
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").Value = "234.39"
$ws.Range("E5").Value = "  -9.22%  "
$ws.Range("D6").Value = "0.602"
$ws.Range("E6").Value = "  -2.89%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "55.01"
$ws.Range("E8").Value = "  -2.21%  "
$ws.Range("D9").Value = "0.371"
$ws.Range("E9").Value = "  -2.15%  "
$ws.Range("D10").Value = "57.53"
$ws.Range("E10").Value = "  +2.53%  "
$ws.Range("D11").Value = "0.0750"
$ws.Range("E11").Value = "  -2.45%  "
$ws.Range("E12").Value = "  -1.66%  "
$ws.Range("D13").Value = "2.319.16"
$ws.Range("E13").Value = "  +0.47%  "
$ws.Range("D14").Value = "14.24"
$ws.Range("E14").Value = "  -0.66%  "
$ws.Range("D15").Value = "20.27"
$ws.Range("E15").Value = "  -5.65%  "
$ws.Range("D16").Value = "0.762"
$ws.Range("E16").Value = "  -4.56%  "
$ws.Range("D17").Value = "5.10"
$ws.Range("E17").Value = "  -2.23%  "
$ws.Range("D18").Value = "2.014.86"
$ws.Range("E18").Value = "  +0.12%  "
$ws.Range("D19").Value = "36.706.51"
$ws.Range("E19").Value = "  -0.98%  "
$ws.Range("D20").Value = "67.78"
$ws.Range("E20").Value = "  -4.64%  "
$ws.Range("D21").Value = "0.0₃0800"
$ws.Range("E21").Value = "  -4.12%  "
$ws.Range("D22").Value = "5.42"
$ws.Range("E22").Value = "  +6.28%  "
$ws.Range("D23").Value = "221.49"
$ws.Range("E23").Value = "  -5.34%  "
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("E25").Value = "  +1.15%  "
$ws.Range("D26").Value = "2.39"
$ws.Range("E26").Value = "  -7.22%  "
$ws.Range("D27").Value = "162.47"
$ws.Range("E27").Value = "  -1.63%  "
$ws.Range("E28").Value = "  +5.85%  "
$ws.Range("D29").Value = "8.64"
$ws.Range("E29").Value = "  -3.53%  "
$ws.Range("D30").Value = "1.38"
$ws.Range("E30").Value = "  +2.34%  "
$ws.Range("D31").Value = "19.03"
$ws.Range("E31").Value = "  -2.67%  "
$ws.Range("D32").Value = "0.117"
$ws.Range("E32").Value = "  -1.51%  "
$ws.Range("D33").Value = "4.38"
$ws.Range("E33").Value = "  -4.69%  "
$ws.Range("D34").Value = "0.0602"
$ws.Range("E34").Value = "  -6.16%  "
$ws.Range("E35").Value = "  +4.20%  "
$ws.Range("D36").Value = "4.26"
$ws.Range("E36").Value = "  -3.60%  "
$ws.Range("E37").Value = "  +0.02%  "
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").Value = "3.34"
$ws.Range("E38").Value = "  -4.63%  "
$ws.Range("B39").Value = "WEMIXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D39").Value = "1.77"
$ws.Range("E39").Value = "  -2.62%  "
$ws.Range("D40").Value = "5.77"
$ws.Range("E40").Value = "  +4.23%  "
$ws.Range("D41").Value = "0.0961"
$ws.Range("E41").Value = "  +4.40%  "
$ws.Range("E42").Value = "  -4.18%  "
$ws.Range("D43").Value = "1.457.86"
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").Value = "0.0204"
$ws.Range("E44").Value = "  -3.01%  "
$ws.Range("B45").Value = "FTXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D45").Value = "4.12"
$ws.Range("E45").Value = "  +38.00%  "
$ws.Range("D46").Value = "90.65"
$ws.Range("E46").Value = "  +1.37%  "
$ws.Range("D47").Value = "1.11"
$ws.Range("E47").Value = "  -6.08%  "
$ws.Range("D48").Value = "15.57"
$ws.Range("E48").Value = "  -0.07%  "
$ws.Range("D49").Value = "1.01"
$ws.Range("E49").Value = "  -1.44%  "
$ws.Range("E50").Value = "  -1.83%  "
$ws.Range("D51").Value = "6.88"
$ws.Range("E51").Value = "  -1.17%  "
